$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("signin")

$ws.Range("C6").Value = '[{"money":"100"}]'
$ws.Range("C7").Value = '[{"money":"300"}]'
$ws.Range("C9").Value = '[{"money":"500"}]'
$ws.Range("C10").Value = '[{"money":"1000"}]'
$ws.Range("C11").Value = '[{"hero":{"id":"2"}}]'
